$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (weekly update), shifting existing data down
$ws.Rows.Item(2).Insert()

# The insert pulls formatting from the header row above; reset it to plain
$ws.Range("A2:T2").ClearFormats()

# Populate the new row 2 with this week's data
$ws.Cells.Item(2, 1).Value = 7
$ws.Cells.Item(2, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(2, 3).Value = "Ñuble"
$ws.Cells.Item(2, 4).Value = 45092
$ws.Cells.Item(2, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(2, 5).Value = 16
$ws.Cells.Item(2, 6).Value = "Fruta"
$ws.Cells.Item(2, 7).Value = 100107
$ws.Cells.Item(2, 8).Value = "Otros"
$ws.Cells.Item(2, 9).Value = 100107001
$ws.Cells.Item(2, 10).Value = "Caqui"
$ws.Cells.Item(2, 11).Value = "Mankaki"
$ws.Cells.Item(2, 12).Value = "Primera"
$ws.Cells.Item(2, 13).Value = 35
$ws.Cells.Item(2, 14).Value = 18000
$ws.Cells.Item(2, 15).Value = 19000
$ws.Cells.Item(2, 16).Value = 18571
$ws.Cells.Item(2, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(2, 18).Value = "Región del Maule"
$ws.Cells.Item(2, 19).Value = 1032
$ws.Cells.Item(2, 20).Value = 18
